$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.150.09"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.656.66"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.53%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.10"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5245"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.38%  "
$ws.Range("E7").Value = "  -0.50%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2623"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06301"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.59"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07804"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.496"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.661.71"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.884.82"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5553"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8018"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.07"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.158.19"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("E19").Value = "  -0.50%  "
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "195.52"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.12"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.960"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.99%  "
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.64"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1204"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.163"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.95"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.494"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05701"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.14%  "
$ws.Range("E31").Value = "  -0.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.485"
$ws.Range("D32").ClearFormats()
$ws.Range("E33").Value = "  +2.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.586"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.42%  "
$ws.Range("B35").Value = "MXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.804"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.57%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9525"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.418"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5701"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01598"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.949"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.061.35"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.40%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8466"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.03%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.005"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "103.52"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.795.37"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.90"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05418"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +5.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.008"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4400"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.42%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₈103"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.24%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.986"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.81%  "
